# Natmi following Dr Hou advice
#
# The underlying NATMI run was redone with an extra "ECs" sending/target
# cluster added to the Slit2 -> App ligand-receptor table, and the
# summary statistics were recomputed (the table now covers the full
# 3x3 grid of ECs / FAPs / sCs senders against ECs / FAPs / sCs
# targets instead of just FAPs / sCs senders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the numeric NATMI edge statistics.
$data = @(
  @("ECs",  "Slit2", "App", "ECs",  2, 0.6666666666666666, 0.01246433333333333, 0.037393,  0.0065371131913745, 0.006537113191374499, 3, 1, 104.794801, 314.384403, 0.3872421191355361, 0.3872421191355361, 1.306197331264333,   11.755775981379,    0.002531445565256729,  0.002531445565256729),
  @("ECs",  "Slit2", "App", "FAPs", 2, 0.6666666666666666, 0.01246433333333333, 0.037393,  0.0065371131913745, 0.006537113191374499, 3, 1, 136.674446, 410.023338, 0.5050451128841343, 0.5050451128841343, 1.703555853092666,   15.332002677834,    0.003301537069674098,  0.003301537069674098),
  @("ECs",  "Slit2", "App", "sCs",  2, 0.6666666666666666, 0.01246433333333333, 0.037393,  0.0065371131913745, 0.006537113191374499, 3, 1, 29.14904533333333, 87.447136,  0.1077127679803296, 0.1077127679803296, 0.3633234173831111,  3.269910756448,     0.0007041305564436733, 0.0007041305564436732),
  @("FAPs", "Slit2", "App", "ECs",  3, 1,                  0.9943730000000001,  2.983119,  0.5215143627507798, 0.5215143627507798,   3, 1, 104.794801, 314.384403, 0.3872421191355361, 0.3872421191355361, 104.205120654773,    937.8460858929571,  0.2019523269912307,    0.2019523269912307),
  @("FAPs", "Slit2", "App", "FAPs", 3, 1,                  0.9943730000000001,  2.983119,  0.5215143627507798, 0.5215143627507798,   3, 1, 136.674446, 410.023338, 0.5050451128841343, 0.5050451128841343, 135.905378892358,    1223.148410031222,  0.2633882802061649,    0.2633882802061649),
  @("FAPs", "Slit2", "App", "sCs",  3, 1,                  0.9943730000000001,  2.983119,  0.5215143627507798, 0.5215143627507798,   3, 1, 29.14904533333333, 87.447136,  0.1077127679803296, 0.1077127679803296, 28.98502365524267,   260.8652128971841,  0.05617375555338417,   0.05617375555338417),
  @("sCs",  "Slit2", "App", "ECs",  3, 1,                  0.8998656666666666,  2.699597,  0.4719485240578458, 0.4719485240578457,   3, 1, 104.794801, 314.384403, 0.3872421191355361, 0.3872421191355361, 94.30124346506567,   848.711191185591,   0.1827583465790487,    0.1827583465790487),
  @("sCs",  "Slit2", "App", "FAPs", 3, 1,                  0.8998656666666666,  2.699597,  0.4719485240578458, 0.4719485240578457,   3, 1, 136.674446, 410.023338, 0.5050451128841343, 0.5050451128841343, 122.9886414660873,   1106.897773194786,  0.2383552956082953,    0.2383552956082953),
  @("sCs",  "Slit2", "App", "sCs",  3, 1,                  0.8998656666666666,  2.699597,  0.4719485240578458, 0.4719485240578457,   3, 1, 29.14904533333333, 87.447136,  0.1077127679803296, 0.1077127679803296, 26.23022511157689,   236.072026004192,   0.05083488187050173,   0.05083488187050172)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}
